$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Stamp the "About" sheet with the date this workbook was brought into the
# new repository (shows as a short date, e.g. 4/21/21).
$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "m/d/yy"
